# "10per change" is the workbook's default ActiveSheet, matching the
# sheet touched by the diff (dimension A1:H52 -> A1:H54).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# E51 / E52 ("bsecode") were stored as text ("543287"/"590024") but should
# become numeric values, matching the rest of the column.
$ws.Range("E51").Value = 543287
$ws.Range("E52").Value = 590024

# Append row 53: a new LODHA (Macrotech Developers Ltd) breakout reading.
$ws.Range("A53").Value = "25/06/2024 08:44:45"
$ws.Range("B53").Value = 1
$ws.Range("C53").Value = "LODHA"
$ws.Range("D53").Value = "Macrotech Developers Ltd"
# bsecode here stays textual (as in the original rows) - force text with a
# leading apostrophe, then strip the resulting quote-prefix style so the
# cell keeps the workbook's default (unstyled) formatting.
$ws.Range("E53").Value = "'543287"
$ws.Range("E53").Style = "Normal"
$ws.Range("F53").Value = -5.19
$ws.Range("G53").Value = 1468.6
$ws.Range("H53").Value = 469410

# Append row 54: a new FACT (Fertilizers And Chemicals Travancore Limited)
# breakout reading.
$ws.Range("A54").Value = "25/06/2024 08:44:45"
$ws.Range("B54").Value = 2
$ws.Range("C54").Value = "FACT"
$ws.Range("D54").Value = "Fertilizers And Chemicals Travancore Limited"
$ws.Range("E54").Value = "'590024"
$ws.Range("E54").Style = "Normal"
$ws.Range("F54").Value = -2.95
$ws.Range("G54").Value = 990.85
$ws.Range("H54").Value = 1755985
